$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.189.88'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '2.242.37'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '293.62'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.59'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.517'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.478'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '31.35'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +7.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.99'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.43'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.78%  '
$ws.Range('D15').Value = '2.582.40'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.27'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').Value = '2.242.39'
$ws.Range('E17').Value = '  +1.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.736'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.53%  '
$ws.Range('D19').Value = '40.080.91'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('D20').Value = '0.0₃0890'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('E21').Value = '  +9.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.86'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.95'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.84'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.73%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.49'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.91%  '
$ws.Range('E27').Value = '  +2.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.01'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('E29').Value = '  +2.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.35'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.45'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.93'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.99'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0726'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.21%  '
$ws.Range('E36').Value = '  +1.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.40'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +7.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.84'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +7.37%  '
$ws.Range('E39').Value = '  +2.04%  '
$ws.Range('E40').Value = '  +3.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.73'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.85'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.55%  '
$ws.Range('D43').Value = '2.062.72'
$ws.Range('E43').Value = '  +7.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.43'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +13.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0270'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.37%  '
$ws.Range('E46').Value = '  +4.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.86'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +8.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.62'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '72.43'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.44%  '
$ws.Range('D50').Value = '2.444.52'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '89.74'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.74%  '
